$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.261.43"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "3.151.74"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.35"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.90"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.148.38"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.30"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "3.672.83"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "3.150.28"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "63.215.45"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.87"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.72"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.67"
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.02"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  -4.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.93"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  -5.56%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.48"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -7.63%  "
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.33"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  -9.17%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.951.28"
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.26"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("E44").Value = "  -5.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").Value = "  -4.90%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.54"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  -8.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.22"
$ws.Range("E51").Value = "  -1.11%  "
